# Add a bold, dark-red "!!" run right after "версий" in the header
# paragraph (before the existing _GoBack bookmark), matching:
#
#   <w:r>
#     <w:rPr>
#       <w:b/>
#       <w:color w:val="C00000"/>
#       <w:sz w:val="40"/>
#       <w:szCs w:val="40"/>
#     </w:rPr>
#     <w:t>!!</w:t>
#   </w:r>
#
# Notably this new run must NOT carry a <w:lang> tag, unlike the
# neighbouring Cyrillic runs. Simply inserting text right after "версий"
# would inherit that run's rPr (including w:lang="ru-RU"), so instead we
# stage the "!!" text at a scratch location that already has no language
# tag (the plain "GIT" run), give it the right Bold/Color/Size there,
# Cut it, and Paste it back in place of a same-length placeholder right
# after "версий" - that keeps its freshly-built (lang-free) formatting
# and lands the new run before the bookmark, exactly like the diff.

$d = $word.ActiveDocument

# 1) Build a correctly-formatted, language-free "!!" run in scratch space
#    at the start of the document (inside the plain "GIT" run: bold,
#    sz=40, no color, no lang), then cut it to the clipboard.
$scratch = $d.Range(3, 3)
$scratch.InsertBefore("!!")
$scratch.Font.Bold = $true
$scratch.Font.Color = 192          ; # 0x0000C0 == RGB(C0,00,00) == "C00000"
$scratch.Font.Size = 20            ; # half-points*2 => w:sz/w:szCs = 40
$scratch.Cut()

# 2) Locate "версий" and collapse to its end (right before the
#    bookmark), then drop in a same-length placeholder there so we have
#    a concrete range to paste over. NOTE: Find.Execute mutates the
#    range it is called on in place, so we must keep operating on the
#    same $target variable rather than re-fetching $d.Content.
$target = $d.Content
$found = $target.Find.Execute("версий", $false, $false, $false, $false,
                               $false, $true, 1, $false, "", 0)
$target.Collapse(0)
$target.InsertBefore("!!")

# 3) Paste the scratch-built, lang-free "!!" run over the placeholder -
#    this preserves the formatting built in step 1 and keeps the new
#    run positioned before the bookmark.
$target.Paste()
